$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "K" column (G) previously held pitch-level "Strike#" counts. This
# regenerates it to hold actual strikeout (K) totals per game, per the
# updated save_data pipeline.
$kValues = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 1
    26 = 2
    27 = 2
    28 = 0
    29 = 3
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 0
    39 = 2
    40 = 0
    41 = 0
    42 = 3
    43 = 1
    44 = 2
    45 = 3
    46 = 2
    47 = 1
    48 = 0
    49 = 2
    50 = 1
    51 = 1
    52 = 1
    53 = 2
    54 = 3
    55 = 1
    56 = 2
    57 = 1
    58 = 1
    60 = 0
    61 = 3
    63 = 2
    64 = 0
    65 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
